$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Yonghui"
$ws.Range("C5").Value = "Anna"
$ws.Range("E4").Value = "Ashish"

$ws.Range("P10").Select()
